$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Phase 1: touch cells in the exact order needed to reproduce
# the shared-string insertion order from the target file.
$ws.Range("E1").Value = "Comments"
$ws.Range("C2").Value = "Medium"
$ws.Range("E2").Value = "With the priority, we are putting in High, medium, and low. We mean a high priority would always be presented, so even if a verbosity setting was set to low, it would still be presented."
$ws.Range("C6").Value = "High"
$ws.Range("C20").Value = "Low"
$ws.Range("E3").Value = "Normally this would be placed on the section element before the H1 or h2 that follows it. The announcement would be better than the generic `"landmark`"."
$ws.Range("E20").Value = "We think that this may appear on each example, which could become too tedious."
$ws.Range("E24").Value = "We felt letting people that the link was to a specific type of element would be good."
$ws.Range("E28").Value = "Notices are many times cautions, posion, warnings, which should always be presented."
$ws.Range("E6").Value = "This would let the reader know they were going back to the referencing item."

# Phase 2: fill in the remaining cells (reusing existing shared strings).
$ws.Range("C3").Value = "Medium"
$ws.Range("C4").Value = "Medium"
$ws.Range("C5").Value = "Medium"
$ws.Range("C7").Value = "Medium"
$ws.Range("C8").Value = "High"
$ws.Range("C9").Value = "Medium"
$ws.Range("C10").Value = "Medium"
$ws.Range("C11").Value = "Medium"
$ws.Range("C12").Value = "Medium"
$ws.Range("C13").Value = "Medium"
$ws.Range("C14").Value = "Medium"
$ws.Range("C15").Value = "Medium"
$ws.Range("C16").Value = "Medium"
$ws.Range("C17").Value = "Medium"
$ws.Range("C18").Value = "Medium"
$ws.Range("C19").Value = "Medium"
$ws.Range("C21").Value = "Medium"
$ws.Range("C22").Value = "Medium"
$ws.Range("C23").Value = "Medium"
$ws.Range("C24").Value = "High"
$ws.Range("C25").Value = "Medium"
$ws.Range("C26").Value = "Medium"
$ws.Range("C27").Value = "High"
$ws.Range("C28").Value = "High"
$ws.Range("C29").Value = "Medium"
$ws.Range("C30").Value = "Medium"
$ws.Range("C31").Value = "Medium"
$ws.Range("C32").Value = "Medium"
$ws.Range("C33").Value = "Medium"
$ws.Range("C34").Value = "Medium"
$ws.Range("C35").Value = "Medium"
$ws.Range("C36").Value = "Medium"
$ws.Range("C37").Value = "Medium"
$ws.Range("C38").Value = "Medium"

# Leave the selection where the author ended up (bottom-right pane).
$ws.Range("C46").Select()
